$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 was right-to-left; the new layout is left-to-right
$ws2.Activate()
$excel.ActiveWindow.DisplayRightToLeft = $false

# Populate Sheet2 data table - values set in an order that reproduces the
# original shared-string table ordering (new strings are interned the
# first time each distinct value is written).
$ws2.Range("A2").Value = "refnum name"
$ws2.Range("B2").Value = "Device"
$ws2.Range("C2").Value = "Function "
$ws2.Range("D2").Value = "Slot"

$ws2.Range("A4").Value = "Drive"
$ws2.Range("B4").Value = "4 Motors"
$ws2.Range("C4").Value = "drive"
$ws2.Range("D4").Value = "PWM 1-4"

$ws2.Range("A5").Value = "Pilot"
$ws2.Range("B5").Value = "Joystick"
$ws2.Range("D5").Value = "USB 1"

$ws2.Range("A6").Value = "CoPilot"
$ws2.Range("B6").Value = "Joystick"
$ws2.Range("D6").Value = "USB 2"

$ws2.Range("A8").Value = "First Cannon"
$ws2.Range("A9").Value = "Second Cannon"
$ws2.Range("A10").Value = "Disc Out"

$ws2.Range("A12").Value = "CannonAI"
$ws2.Range("A13").Value = "Enc"

$ws2.Range("A15").Value = "Flipper"
$ws2.Range("A16").Value = "Conveyer Belt"
$ws2.Range("A17").Value = "Roller"

$ws2.Range("A19").Value = "Cannon"
$ws2.Range("A20").Value = "InFlipper"
$ws2.Range("A22").Value = "Middle"
$ws2.Range("A21").Value = "FlipperAtPlace"
$ws2.Range("A23").Value = "Bottom"
$ws2.Range("A24").Value = "FrisbeeDirection"

$ws2.Range("C5").Value = "driver's joystick"
$ws2.Range("C6").Value = "operator's joystick"

# Column widths (character units). The host's column-width pixel grid
# snaps to the nearest 1/7 character, so these inputs are chosen to land
# on the closest achievable value to the authored widths (14.375, 17.5,
# 15.5).
$ws2.Columns.Item(1).ColumnWidth = 13.714285714285714
$ws2.Columns.Item(2).ColumnWidth = 16.714285714285715
$ws2.Columns.Item(3).ColumnWidth = 14.714285714285715

$ws2.Range("B8").Select()
$ws1.Range("B9").Select()

$ws2.Activate()
